$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.026992
$ws.Range("H2").Value = 0.08097599999999999
$ws.Range("I2").Value = 0.004182906599909731
$ws.Range("J2").Value = 0.00420788870005516
$ws.Range("M2").Value = 0.3615393333333333
$ws.Range("N2").Value = 1.084618
$ws.Range("O2").Value = 0.7649240942154193
$ws.Range("P2").Value = 0.7664955283791567
$ws.Range("Q2").Value = 0.009758669685333331
$ws.Range("R2").Value = 0.08782802716799998
$ws.Range("S2").Value = 0.00319960604212365
$ws.Range("T2").Value = 0.003225327872509463

$ws.Range("G3").Value = 0.026992
$ws.Range("H3").Value = 0.08097599999999999
$ws.Range("I3").Value = 0.004182906599909731
$ws.Range("J3").Value = 0.00420788870005516
$ws.Range("O3").Value = 0.228925442648571
$ws.Range("P3").Value = 0.2293957393280025
$ws.Range("Q3").Value = 0.002920561392
$ws.Range("R3").Value = 0.02628505252799999
$ws.Range("S3").Value = 0.0009575737449419641
$ws.Range("T3").Value = 0.0009652717393591009

$ws.Range("D4").Value = "MuSCs"
$ws.Range("G4").Value = 0.026992
$ws.Range("H4").Value = 0.08097599999999999
$ws.Range("I4").Value = 0.004182906599909731
$ws.Range("J4").Value = 0.00420788870005516
$ws.Range("L4").Value = 0.5
$ws.Range("M4").Value = 0.002907
$ws.Range("N4").Value = 0.005814
$ws.Range("O4").Value = 0.006150463136009796
$ws.Range("P4").Value = 0.004108732292840814
$ws.Range("Q4").Value = 0.000078465744
$ws.Range("R4").Value = 0.000470794464
$ws.Range("S4").Value = 0.00002572681284411687
$ws.Range("T4").Value = 0.00001728908818659659

$ws.Range("I5").Value = 0.976381346197431
$ws.Range("J5").Value = 0.9822127115383066
$ws.Range("M5").Value = 0.3615393333333333
$ws.Range("N5").Value = 1.084618
$ws.Range("O5").Value = 0.7649240942154193
$ws.Range("P5").Value = 0.7664955283791567
$ws.Range("Q5").Value = 2.277885679940222
$ws.Range("R5").Value = 20.500971119462
$ws.Range("S5").Value = 0.7468576168489016
$ws.Range("T5").Value = 0.7528616513112785

$ws.Range("I6").Value = 0.976381346197431
$ws.Range("J6").Value = 0.9822127115383066
$ws.Range("O6").Value = 0.228925442648571
$ws.Range("P6").Value = 0.2293957393280025
$ws.Range("S6").Value = 0.2235185318720545
$ws.Range("T6").Value = 0.2253154111406919

$ws.Range("D7").Value = "MuSCs"
$ws.Range("I7").Value = 0.976381346197431
$ws.Range("J7").Value = 0.9822127115383066
$ws.Range("L7").Value = 0.5
$ws.Range("M7").Value = 0.002907
$ws.Range("N7").Value = 0.005814
$ws.Range("O7").Value = 0.006150463136009796
$ws.Range("P7").Value = 0.004108732292840814
$ws.Range("Q7").Value = 0.018315610671
$ws.Range("R7").Value = 0.109893664026
$ws.Range("S7").Value = 0.006005197476474917
$ws.Range("T7").Value = 0.004035649086336179

$ws.Range("G8").Value = 0.1149325
$ws.Range("H8").Value = 0.229865
$ws.Range("I8").Value = 0.01781090370458377
$ws.Range("J8").Value = 0.01194485200600399
$ws.Range("M8").Value = 0.3615393333333333
$ws.Range("N8").Value = 1.084618
$ws.Range("O8").Value = 0.7649240942154193
$ws.Range("P8").Value = 0.7664955283791567
$ws.Range("Q8").Value = 0.04155261942833333
$ws.Range("R8").Value = 0.24931571657
$ws.Range("S8").Value = 0.0136239893833868
$ws.Range("T8").Value = 0.009155675649752861

$ws.Range("G9").Value = 0.1149325
$ws.Range("H9").Value = 0.229865
$ws.Range("I9").Value = 0.01781090370458377
$ws.Range("J9").Value = 0.01194485200600399
$ws.Range("O9").Value = 0.228925442648571
$ws.Range("P9").Value = 0.2293957393280025
$ws.Range("Q9").Value = 0.0124358114325
$ws.Range("R9").Value = 0.07461486859499999
$ws.Range("S9").Value = 0.004077369014542912
$ws.Range("T9").Value = 0.00274009815708086

$ws.Range("D10").Value = "MuSCs"
$ws.Range("G10").Value = 0.1149325
$ws.Range("H10").Value = 0.229865
$ws.Range("I10").Value = 0.01781090370458377
$ws.Range("J10").Value = 0.01194485200600399
$ws.Range("L10").Value = 0.5
$ws.Range("M10").Value = 0.002907
$ws.Range("N10").Value = 0.005814
$ws.Range("O10").Value = 0.006150463136009796
$ws.Range("P10").Value = 0.004108732292840814
$ws.Range("Q10").Value = 0.0003341087775
$ws.Range("R10").Value = 0.00133643511
$ws.Range("S10").Value = 0.0001095453066540628
$ws.Range("T10").Value = 0.00004907819917027298

$ws.Range("G11").Value = 0.010485
$ws.Range("H11").Value = 0.031455
$ws.Range("I11").Value = 0.001624843498075486
$ws.Range("J11").Value = 0.001634547755634201
$ws.Range("M11").Value = 0.3615393333333333
$ws.Range("N11").Value = 1.084618
$ws.Range("O11").Value = 0.7649240942154193
$ws.Range("P11").Value = 0.7664955283791567
$ws.Range("Q11").Value = 0.003790739909999999
$ws.Range("R11").Value = 0.03411665919
$ws.Range("S11").Value = 0.001242881941007205
$ws.Range("T11").Value = 0.001252873545615802

$ws.Range("G12").Value = 0.010485
$ws.Range("H12").Value = 0.031455
$ws.Range("I12").Value = 0.001624843498075486
$ws.Range("J12").Value = 0.001634547755634201
$ws.Range("O12").Value = 0.228925442648571
$ws.Range("P12").Value = 0.2293957393280025
$ws.Range("Q12").Value = 0.001134487485
$ws.Range("R12").Value = 0.010210387365
$ws.Range("S12").Value = 0.0003719680170315832
$ws.Range("T12").Value = 0.0003749582908706347

$ws.Range("D13").Value = "MuSCs"
$ws.Range("G13").Value = 0.010485
$ws.Range("H13").Value = 0.031455
$ws.Range("I13").Value = 0.001624843498075486
$ws.Range("J13").Value = 0.001634547755634201
$ws.Range("L13").Value = 0.5
$ws.Range("M13").Value = 0.002907
$ws.Range("N13").Value = 0.005814
$ws.Range("O13").Value = 0.006150463136009796
$ws.Range("P13").Value = 0.004108732292840814
$ws.Range("Q13").Value = 0.000030479895
$ws.Range("R13").Value = 0.00018287937
$ws.Range("S13").Value = 0.000009993540036698483
$ws.Range("T13").Value = 0.000006715919147764716
